# Import a Flow excel file should be idempotent (close #53)
#
# The "Message_Flow" sheet gains a new "External" column (inserted right
# after "Alias flow", i.e. before the old "Source Element" column) and a
# new data row describing an external flow (EXT.001).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Message_Flow")
$ws.Activate()

# --- Insert the new "External" column before the current column C ---------
$ws.Columns("C:C").Insert()

# --- Grow the table (ListObject) so it still starts at the data header ----
# A plain column insert does not itself re-anchor the table/autofilter ref,
# so resize it explicitly to the now-shifted D1:U1 header range. This has
# to happen *before* we populate the new C1 header cell, otherwise the
# table (whose range still covers the new column at this point) picks up
# the new cell text as its first column's name.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D1:U1"))

# Header for the newly inserted column.
$ws.Range("C1").Value2 = "External"

# Keep the narrow width the new column has in the target workbook.
$ws.Columns("C:C").ColumnWidth = 7.99

# --- Append the new row describing the external flow -----------------------
$ws.Range("A6").Value2 = "EXT.001"
$ws.Range("B6").Value2 = "S.02"
$ws.Range("C6").Value2 = "yes"
$ws.Range("D6").Value2 = "APPLICATION-0004"
$ws.Range("E6").Value2 = "APPLICATION-0003"
$ws.Range("F6").Value2 = "Description TRAD.004"
$ws.Range("G6").Value2 = "EXTERNAL"
$ws.Range("H6").Value2 = "OTHER"
$ws.Range("J6").Value2 = "JSON"

# Match the workbook's final selection/cursor position.
$ws.Range("B6").Select()
